# Update the "Native" results (columns B:C, rows 7-16) with the latest
# benchmark run data, update the Native Base Score (B18), and leave the
# selection on the Native results header (B17:C17) as the last editing
# focus, matching what was selected when the data was refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Native "Seconds" (B) / "Ratio" (C) per benchmark row
$ws.Range("B7").Value  = 654
$ws.Range("C7").Value  = 2.71

$ws.Range("B8").Value  = 953
$ws.Range("C8").Value  = 4.18

$ws.Range("B9").Value  = 1409
$ws.Range("C9").Value  = 3.35

$ws.Range("B10").Value = 779
$ws.Range("C10").Value = 2.09

$ws.Range("B11").Value = 533
$ws.Range("C11").Value = 2.66

$ws.Range("B12").Value = 554
$ws.Range("C12").Value = 3.19

$ws.Range("B13").Value = 773
$ws.Range("C13").Value = 1.85

$ws.Range("B14").Value = 1045
$ws.Range("C14").Value = 1.63

$ws.Range("B15").Value = 780
$ws.Range("C15").Value = 3.77

$ws.Range("B16").Value = 4665
$ws.Range("C16").Value = 1.33

# Native Base Score
$ws.Range("B18").Value = 2.52

# Leave selection on the Native header cell, like after editing this block
$ws.Range("B17:C17").Select()
